# Commit: "Change date in rebate2 file"
# The only content change in the target diff is the date string stored in
# cell G7 ("End User Invoice Date") — 08/27/2018 -> 11/25/2018. Update it
# through the Range.Value API so the workbook's shared-string table picks
# up the new text in place (same cell / same shared-string slot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = "11/25/2018"
